# Auto-applies the cell-value updates described by the commit diff.
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H54").Value = 15000
$ws.Range("I54").Value = 10000
$ws.Range("J54").Value = 20000
$ws.Range("K54").Value = 10000
$ws.Range("L54").Value = 20000
$ws.Range("M54").Value = -9514
$ws.Range("N54").Value = -20972

$ws.Range("H70").Value = 4593.8
$ws.Range("I70").Value = 2359.4358
$ws.Range("J70").Value = 19117.166
$ws.Range("K70").Value = 7078.307400000001
$ws.Range("L70").Value = 57351.49800000001
$ws.Range("M70").Value = -6808.307400000001
$ws.Range("N70").Value = -57891.49800000001

$ws.Range("H73").Value = 4593.8
$ws.Range("I73").Value = 2359.4358
$ws.Range("J73").Value = 19117.166
$ws.Range("K73").Value = 7078.307400000001
$ws.Range("L73").Value = 57351.49800000001
$ws.Range("M73").Value = -6142.307400000001
$ws.Range("N73").Value = -59223.49800000001

$ws.Range("H74").Value = 5853.48
$ws.Range("I74").Value = 5253.5
$ws.Range("K74").Value = 5253.5
$ws.Range("M74").Value = -4317.5

$ws.Range("H77").Value = 5853.48
$ws.Range("I77").Value = 5253.5
$ws.Range("K77").Value = 26267.5
$ws.Range("M77").Value = -21587.5

$ws.Range("H97").Value = 7323
$ws.Range("J97").Value = 7323
$ws.Range("L97").Value = 21969
$ws.Range("N97").Value = -22961

$ws.Range("H100").Value = 7346.3076
$ws.Range("I100").Value = 2924.8333
$ws.Range("J100").Value = 11136.143
$ws.Range("K100").Value = 2924.8333
$ws.Range("L100").Value = 11136.143
$ws.Range("M100").Value = -2383.8333
$ws.Range("N100").Value = -12218.143

$ws.Range("H103").Value = 792.9375
$ws.Range("I103").Value = 815.5714
$ws.Range("J103").Value = 775.3333
$ws.Range("K103").Value = 2446.7142
$ws.Range("L103").Value = 2325.9999
$ws.Range("M103").Value = -1860.7142
$ws.Range("N103").Value = -3497.9999

$ws.Range("H111").Value = 2961.0908
$ws.Range("I111").Value = 2600
$ws.Range("J111").Value = 3262
$ws.Range("K111").Value = 7800
$ws.Range("L111").Value = 9786
$ws.Range("M111").Value = -4733
$ws.Range("N111").Value = -15920

$ws.Range("H112").Value = 2093.9
$ws.Range("I112").Value = 664.2
$ws.Range("J112").Value = 2570.4666
$ws.Range("K112").Value = 1992.6
$ws.Range("L112").Value = 7711.399800000001
$ws.Range("M112").Value = -884.6000000000001
$ws.Range("N112").Value = -9927.399800000001

$ws.Range("H125").Value = 1941.6154
$ws.Range("I125").Value = 2007
$ws.Range("J125").Value = 1885.5714
$ws.Range("K125").Value = 18063
$ws.Range("L125").Value = 16970.1426
$ws.Range("M125").Value = -15603
$ws.Range("N125").Value = -21890.1426

$ws.Range("H132").Value = 3166.6
$ws.Range("I132").Value = 2993.5117
$ws.Range("K132").Value = 8980.535100000001
$ws.Range("M132").Value = -6450.535100000001

$ws.Range("H137").Value = 1573.0555
$ws.Range("I137").Value = 1422.8667
$ws.Range("J137").Value = 2324
$ws.Range("K137").Value = 4268.6001
$ws.Range("L137").Value = 6972
$ws.Range("M137").Value = -1718.6001
$ws.Range("N137").Value = -12072

$ws.Range("H141").Value = 7279.9565
$ws.Range("I141").Value = 6022
$ws.Range("J141").Value = 15666.333
$ws.Range("K141").Value = 18066
$ws.Range("L141").Value = 46998.999
$ws.Range("M141").Value = -12886
$ws.Range("N141").Value = -57358.999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2119.3809
$ws.Range("I32").Value = 1342.45
$ws.Range("K32").Value = 1342.45
$ws.Range("M32").Value = -1055.45

$ws.Range("H45").Value = 3507.1
$ws.Range("I45").Value = 4399.3335
$ws.Range("K45").Value = 4399.3335
$ws.Range("M45").Value = -4022.3335

$ws.Range("H132").Value = 2152.52
$ws.Range("J132").Value = 2715
$ws.Range("L132").Value = 8145
$ws.Range("N132").Value = -13205

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2254.238
$ws.Range("I107").Value = 1551.8125
$ws.Range("K107").Value = 1551.8125
$ws.Range("M107").Value = 368.1875

$ws.Range("H132").Value = 95999
$ws.Range("J132").Value = 95999
$ws.Range("L132").Value = 95999
$ws.Range("N132").Value = -106119

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 736
$ws.Range("I16").Value = 768.625
$ws.Range("J16").Value = 475
$ws.Range("K16").Value = 768.625
$ws.Range("L16").Value = 475
$ws.Range("M16").Value = -481.625
$ws.Range("N16").Value = -1049

$ws.Range("H31").Value = 3098
$ws.Range("I31").Value = 2859.8
$ws.Range("K31").Value = 2859.8
$ws.Range("M31").Value = -2564.8

$ws.Range("H34").Value = 3098
$ws.Range("I34").Value = 2859.8
$ws.Range("K34").Value = 2859.8
$ws.Range("M34").Value = -2657.8

$ws.Range("H113").Value = 736
$ws.Range("I113").Value = 768.625
$ws.Range("J113").Value = 475
$ws.Range("K113").Value = 768.625
$ws.Range("L113").Value = 475
$ws.Range("M113").Value = 1401.375
$ws.Range("N113").Value = -4815

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 2134.3333
$ws.Range("I122").Value = 251.28572
$ws.Range("K122").Value = 2261.57148
$ws.Range("M122").Value = 188.4285199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 2100999.8
$ws.Range("J14").Value = 1666.3334
$ws.Range("L14").Value = 1666.3334
$ws.Range("N14").Value = -2002.3334

$ws.Range("H102").Value = 9469.941999999999
$ws.Range("I102").Value = 15032.056
$ws.Range("J102").Value = 3580.647
$ws.Range("K102").Value = 15032.056
$ws.Range("L102").Value = 3580.647
$ws.Range("M102").Value = -13410.056
$ws.Range("N102").Value = -6824.647

$ws.Range("H113").Value = 3512.5652
$ws.Range("I113").Value = 2518
$ws.Range("J113").Value = 3947.6875
$ws.Range("K113").Value = 2518
$ws.Range("L113").Value = 3947.6875
$ws.Range("M113").Value = -348
$ws.Range("N113").Value = -8287.6875

$ws.Range("H122").Value = 3535.6333
$ws.Range("I122").Value = 3017.8948
$ws.Range("J122").Value = 4429.909
$ws.Range("K122").Value = 9053.6844
$ws.Range("L122").Value = 13289.727
$ws.Range("M122").Value = -6603.6844
$ws.Range("N122").Value = -18189.727

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1824.6666
$ws.Range("I61").Value = 1586.5555
$ws.Range("K61").Value = 1586.5555
$ws.Range("M61").Value = -1384.5555

$ws.Range("H113").Value = 1824.6666
$ws.Range("I113").Value = 1586.5555
$ws.Range("K113").Value = 1586.5555
$ws.Range("M113").Value = 583.4445000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 34702.8
$ws.Range("J51").Value = 33611
$ws.Range("L51").Value = 33611
$ws.Range("N51").Value = -34631

$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

$ws.Range("H54").Value = 29600.334
$ws.Range("J54").Value = 29600.334
$ws.Range("L54").Value = 29600.334
$ws.Range("N54").Value = -30640.334

$ws.Range("H81").Value = 4995.3335
$ws.Range("I81").Value = 6441.769
$ws.Range("J81").Value = 1234.6
$ws.Range("K81").Value = 12883.538
$ws.Range("L81").Value = 2469.2
$ws.Range("M81").Value = -11822.538
$ws.Range("N81").Value = -4591.2

$ws.Range("H84").Value = 4995.3335
$ws.Range("I84").Value = 6441.769
$ws.Range("J84").Value = 1234.6
$ws.Range("K84").Value = 64417.69
$ws.Range("L84").Value = 12346
$ws.Range("M84").Value = -59113.69
$ws.Range("N84").Value = -22954

$ws.Range("H122").Value = 3025.5
$ws.Range("I122").Value = 2570.9375
$ws.Range("J122").Value = 4843.75
$ws.Range("K122").Value = 7712.8125
$ws.Range("L122").Value = 14531.25
$ws.Range("M122").Value = -5262.8125
$ws.Range("N122").Value = -19431.25

$ws.Range("H126").Value = 2998
$ws.Range("I126").Value = 2999.2
$ws.Range("J126").Value = 2995
$ws.Range("K126").Value = 8997.599999999999
$ws.Range("L126").Value = 8985
$ws.Range("M126").Value = -6527.599999999999
$ws.Range("N126").Value = -13925

$ws.Range("H136").Value = 2219.2258
$ws.Range("I136").Value = 2432.2144
$ws.Range("J136").Value = 2043.8235
$ws.Range("K136").Value = 7296.6432
$ws.Range("L136").Value = 6131.470499999999
$ws.Range("M136").Value = -4746.6432
$ws.Range("N136").Value = -11231.4705
